$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Force D:E range to text so numeric-looking strings (e.g. "1.01", "22.10")
# are preserved verbatim instead of being parsed into floating point numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "27.319.28"
$ws.Range("E2").Value = "  +1.37%  "
$ws.Range("D3").Value = "1.568.87"
$ws.Range("E3").Value = "  +0.41%  "
$ws.Range("E4").Value = "  +0.38%  "
$ws.Range("E5").Value = "  +1.70%  "
$ws.Range("E6").Value = "  +0.56%  "
$ws.Range("D7").Value = "1.01"
$ws.Range("E7").Value = "  +0.28%  "
$ws.Range("D8").Value = "22.10"
$ws.Range("E8").Value = "  +0.19%  "
$ws.Range("D9").Value = "0.249"
$ws.Range("E9").Value = "  +0.55%  "
$ws.Range("E10").Value = "  -0.29%  "
$ws.Range("D11").Value = "0.0871"
$ws.Range("E11").Value = "  +1.62%  "
$ws.Range("D12").Value = "1.792.61"
$ws.Range("E12").Value = "  +0.44%  "
$ws.Range("D13").Value = "1.575.31"
$ws.Range("E13").Value = "  +0.76%  "
$ws.Range("E14").Value = "  +0.83%  "
$ws.Range("E15").Value = "  -0.06%  "
$ws.Range("D16").Value = "27.281.03"
$ws.Range("E16").Value = "  +1.20%  "
$ws.Range("D17").Value = "62.24"
$ws.Range("E17").Value = "  +0.29%  "
$ws.Range("D18").Value = "7.51"
$ws.Range("E18").Value = "  +2.17%  "
$ws.Range("D19").Value = "217.32"
$ws.Range("E19").Value = "  +0.33%  "
$ws.Range("D20").Value = "0.0₃0705"
$ws.Range("E20").Value = "  -0.16%  "
$ws.Range("E21").Value = "  +0.19%  "
$ws.Range("E22").Value = "  +1.21%  "
$ws.Range("D23").Value = "9.23"
$ws.Range("E23").Value = "  +0.23%  "
$ws.Range("E24").Value = "  +0.29%  "
$ws.Range("D25").Value = "153.65"
$ws.Range("E25").Value = "  +0.72%  "
$ws.Range("E26").Value = "  +0.57%  "
$ws.Range("D27").Value = "15.05"
$ws.Range("E27").Value = "  -0.15%  "
$ws.Range("E28").Value = "  +1.85%  "
$ws.Range("E29").Value = "  +0.12%  "
$ws.Range("D30").Value = "1.15"
$ws.Range("E30").Value = "  +2.91%  "
$ws.Range("E31").Value = "  +0.21%  "
$ws.Range("E32").Value = "  +0.50%  "
$ws.Range("E33").Value = "  +1.76%  "
$ws.Range("D34").Value = "1.445.39"
$ws.Range("E34").Value = "  +1.70%  "
$ws.Range("D35").Value = "1.11"
$ws.Range("E35").Value = "  +4.38%  "
$ws.Range("E36").Value = "  -0.11%  "
$ws.Range("E37").Value = "  +0.46%  "
$ws.Range("E38").Value = "  +0.97%  "
$ws.Range("E39").Value = "  +0.07%  "
$ws.Range("D40").Value = "5.92"
$ws.Range("E40").Value = "  +2.25%  "
$ws.Range("D41").Value = "0.810"
$ws.Range("E41").Value = "  +0.26%  "
$ws.Range("E42").Value = "  +0.24%  "
$ws.Range("E43").Value = "  +0.78%  "
$ws.Range("E44").Value = "  -0.33%  "
$ws.Range("D45").Value = "64.66"
$ws.Range("E46").Value = "  -0.81%  "
$ws.Range("D47").Value = "1.704.82"
$ws.Range("E47").Value = "  +0.35%  "
$ws.Range("D48").Value = "85.99"
$ws.Range("E48").Value = "  -1.53%  "
$ws.Range("E49").Value = "  +1.23%  "
$ws.Range("E50").Value = "  +0.04%  "
$ws.Range("D51").Value = "0.0957"
$ws.Range("E51").Value = "  -0.27%  "

# Strip the temporary text-number-format so the cell style matches the
# original (no explicit "s" style index), leaving only the value as text.
$ws.Range("D2:E51").Style = "Normal"
